$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Station"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "DrogDepth"
$ws.Range("D1").Value = "DepDate"
$ws.Range("E1").Value = "DepTime"
$ws.Range("F1").Value = "DepLong"
$ws.Range("G1").Value = "DepLat"
$ws.Range("H1").Value = "RecovDate"
$ws.Range("I1").Value = "RecovTime"
$ws.Range("J1").Value = "RecovLong"
$ws.Range("K1").Value = "RecovLat"

$ws.Range("L1").Select() | Out-Null
